$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as text,
# even when the value looks numeric (e.g. "56.0000") and the cell's number
# format is actually numeric ("0.00"). Excel would otherwise silently
# coerce such strings into real numbers, which would not match the source
# report (every data cell there is stored as a shared string). Cells whose
# number format is already text ("@") don't need this - and re-applying
# "@" to them can swap in a different (but equivalent) style record, so we
# only use this helper where the underlying format truly is numeric.
function Set-TextValue($range, $text) {
    $orig = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = $orig
}

# The "SILDOCARE 8MG 30CAPS" row was removed from the report and a new
# "STERONATE  5MG 20 TAB" row was appended to this alphabetical block.
# That shifts "SPASMO-DIGESTIN 30 TABS." and "SPINOBAC 5MG/5ML SYRUP 120ML"
# up by one row (to rows 50 and 51), and the new row lands at row 52.

# Row 50: was SILDOCARE, now SPASMO-DIGESTIN
$ws.Range("C50").Value = "SPASMO-DIGESTIN 30 TABS."
$ws.Range("H50").Value = "3:2"
$ws.Range("N50").Value = "78.00"
Set-TextValue $ws.Range("P50") "25.7400"
$ws.Range("Q50").Value = "0:1"

# Row 51: was SPASMO-DIGESTIN, now SPINOBAC
$ws.Range("C51").Value = "SPINOBAC 5MG/5ML SYRUP 120ML"
$ws.Range("H51").Value = "0:0"
$ws.Range("N51").Value = "39.00"
Set-TextValue $ws.Range("P51") "39.0000"
$ws.Range("Q51").Value = "1:0"

# Row 52: was SPINOBAC, now the newly added STERONATE row
$ws.Range("C52").Value = "STERONATE  5MG 20 TAB"
$ws.Range("H52").Value = "0:1"
$ws.Range("N52").Value = "56.00"
Set-TextValue $ws.Range("P52") "56.0000"
$ws.Range("Q52").Value = "1:0"

# Grand total (column P) reflects the price swap: -44.5500 (SILDOCARE) + 56.0000 (STERONATE)
$ws.Range("P72").Value = 3244.96

# Footer timestamp updated to reflect the new export time
$ws.Range("A73").Value = "Monday, 16 June, 2025 8:02 PM"
